$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for all data rows (2-36)
# from 45666 (2025-01-09) to 45667 (2025-01-10).
for ($row = 2; $row -le 36; $row++) {
    $ws.Cells.Item($row, 3).Value = 45667
}

# Rows 35 and 36 swapped their "Beteckning" (A) and "Area (ha)" (G) values.
$ws.Range("A35").Value = "A 60500-2024"
$ws.Range("G35").Value = 0.8

$ws.Range("A36").Value = "A 60501-2024"
$ws.Range("G36").Value = 0.6
